$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The daily work report grows by one more day's entries (rows 97-104),
# following the same 7-data-row + 1-blank-separator-row pattern used
# throughout the sheet. Copy formatting from an existing same-shaped
# block (rows 54-60 have the same 3-date-style + 4-plain-style layout,
# row 52 is a blank separator row) so borders/number formats/fonts match.
$ws.Range("A54:D60").Copy() | Out-Null
$ws.Range("A97:D103").PasteSpecial(-4122) | Out-Null

$ws.Range("A52:D52").Copy() | Out-Null
$ws.Range("A104:D104").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Row 97: new day header
$ws.Range("A97").Value = 45686
$ws.Range("B97").Value = "Domm"
$ws.Range("D97").Value = 0.25

# Row 98: Meeting / Reconsile
$ws.Range("B98").Value = "Meeting"
$ws.Range("C98").Value = "Reconsile"
$ws.Range("D98").Value = 1

# Row 99: Light issue + Back account opening
$ws.Range("C99").Value = "Light issue + Back account opening"
$ws.Range("D99").Value = 1.5

# Row 100: General Discussion
$ws.Range("C100").Value = "General Discussion"
$ws.Range("D100").Value = 0.25

# Row 101: Study / Dependency Injection
$ws.Range("B101").Value = "Study"
$ws.Range("C101").Value = "Dependency Injection"
$ws.Range("D101").Value = 3

# Row 102: Reconcile Corrections
$ws.Range("C102").Value = "Reconcile Corrections"
$ws.Range("D102").Value = 1

# Row 103: Reconcile Revision
$ws.Range("C103").Value = "Reconcile Revision"
$ws.Range("D103").Value = 1

# Row 104 stays blank (separator row), matching the pattern.

# Update the visible selection to the newly-added block, as Excel would
# leave it after the author typed this in.
$ws.Range("A97:D104").Select() | Out-Null
